$d = $word.ActiveDocument

# Map of old text -> new text. Each old value is unique in the document,
# so MatchWholeWord keeps the replacements scoped to the right run.
$replacements = [ordered]@{
    "2025-06-06 Friday" = "2025-06-07 Saturday"
    "14×68="            = "86×62="
    "25×33="            = "31×61="
    "28×94="            = "79×22="
    "33×25="            = "59×92="
    "40×27="            = "49×58="
    "89×68="            = "27×21="
    "90×27="            = "40×78="
    "46×74="            = "45×23="
    "33×80="            = "28×87="
    "42×63="            = "92×33="
    "12×37="            = "30×55="
    "95×29="            = "58×46="
    "74×32="            = "74×23="
    "27×29="            = "67×31="
    "27×34="            = "85×68="
    "82×92="            = "84×98="
    "87×85="            = "36×77="
    "34×99="            = "82×45="
    "74×89="            = "99×84="
    "72×84="            = "96×61="
    "81×59="            = "94×61="
    "99×53="            = "57×40="
    "56×57="            = "59×80="
    "87×19="            = "77×14="
    "90×54="            = "43×92="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
